$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column I (Max Port, Ready Port),
# shifting the existing "First Updated"/"Last Updated" columns to K/L.
$ws.Range("I1:J1").EntireColumn.Insert()

# New header labels
$ws.Range("I1").Value = "Max Port"
$ws.Range("J1").Value = "Ready Port"

# New "Max Port" values for rows 2-9, and touch the new "Ready Port" column
# so the (blank) cells get materialized, just like the original sheet.
$ws.Range("I2").Value = 8
$ws.Range("I3").Value = 8
$ws.Range("I4").Value = 16
$ws.Range("I5").Value = 8
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 8
$ws.Range("I8").Value = 16
$ws.Range("I9").Value = 8

# Materialize the new (blank) "Ready Port" cells as empty text cells,
# matching the empty inline-string cells used throughout this sheet.
$ws.Range("J2:J9").Value = "'"
$ws.Range("J2:J9").ClearFormats()

# Data corrections
$ws.Range("G4").Value = "'22990459"
$ws.Range("G4").ClearFormats()
$ws.Range("L4").Value = 45333.70420138889
